# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) values were recomputed and need to be updated
# in place on the active worksheet. Only the G column (rows 2-13) changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") keyed by row number
$kValues = @{
    2  = 3
    3  = 0
    4  = 1
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 2
    13 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
